# Automatische test-sync: 2025-08-26 20:37:50
#
# Adds a new "Retour status" log entry (row 7) to the "Logs" sheet,
# extends the conditional formatting ranges to include the new row,
# and updates the "Retour / Terugbetaling" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A7").Value = "Retour status"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("D7").Value = "Retour / Terugbetaling"
$logs.Range("F7").Value = "2025-08-26 20:37:36"
$logs.Range("G7").Value = "Ja"
$logs.Range("H7").Value = "Nee"
$logs.Range("I7").Value = "Nee"
$logs.Range("J7").Value = "Nee"

# --- Extend the conditional formatting sqref ranges to row 7 -----------
$logs.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D7"))
$logs.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G7"))
$logs.Range("H2:H6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H7"))
$logs.Range("I2:I6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I7"))
$logs.Range("J2:J6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J7"))

# --- Dashboard sheet: bump the "Retour / Terugbetaling" tally ----------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 2
